# Refresh cryptos list snapshot (prices + 1h volume deltas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.879.53'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '2.399.63'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').Value = '565.38'
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.10'
$ws.Range('E6').Value = '  +2.66%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('E8').Value = '  +2.50%  '
$ws.Range('D9').Value = '2.403.80'
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('E10').Value = '  +1.70%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '5.19'
$ws.Range('E12').Value = '  +2.91%  '
$ws.Range('D13').Value = '0.346'
$ws.Range('E13').Value = '  +2.98%  '
$ws.Range('D14').Value = '26.42'
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000170'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '2.835.02'
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('D17').Value = '60.726.30'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '2.412.77'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = '8.05'
$ws.Range('E19').Value = '  +3.62%  '
$ws.Range('D20').Value = '10.71'
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('D21').Value = '324.54'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('E22').Value = '  +1.73%  '
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.90'
$ws.Range('E25').Value = '  +4.42%  '
$ws.Range('D26').Value = '65.07'
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('D27').Value = '586.07'
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('D28').Value = '8.21'
$ws.Range('E28').Value = '  -0.73%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0943'
$ws.Range('E29').Value = '  +2.21%  '
$ws.Range('B30').Value = 'WrappedeETH'
$ws.Range('C30').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D30').Value = '2.509.21'
$ws.Range('E30').Value = '  -1.09%  '
$ws.Range('D31').Value = '8.04'
$ws.Range('E31').Value = '  +2.79%  '
$ws.Range('D32').Value = '1.35'
$ws.Range('E32').Value = '  +0.84%  '
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('E34').Value = '  +0.80%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '1.47'
$ws.Range('E35').Value = '  +5.44%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').Value = '0.996'
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('E37').Value = '  +1.67%  '
$ws.Range('D38').Value = '4.62'
$ws.Range('E38').Value = '  +0.92%  '
$ws.Range('D39').Value = '151.87'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').Value = '18.33'
$ws.Range('E40').Value = '  +0.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.20'
$ws.Range('E41').Value = '  +2.22%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').Value = '2.53'
$ws.Range('E43').Value = '  +9.74%  '
$ws.Range('D44').Value = '1.69'
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('D45').Value = '41.58'
$ws.Range('D46').Value = '0.0₆0279'
$ws.Range('E46').Value = '  +6.48%  '
$ws.Range('D47').Value = '141.81'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').Value = '3.52'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.590'
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0510'
$ws.Range('E50').Value = '  +1.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.50'
$ws.Range('E51').Value = '  +1.27%  '
